$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Total Number of Units"
$ws.Range("C1").Value = "Notes"
$ws.Range("B2").Value = "b"
$ws.Range("C2").Value = "a"
$ws.Range("B3").Value = "d"
$ws.Range("C3").Value = "c"

$ws.Range("F9").Select()
